$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-5 (take former rows 6-9 values)
$newRows = @{
    2 = @{ D = 44890; M = 150; N = 13000; O = 13000; P = 13000; R = "La Ligua";            S = 2600 }
    3 = @{ D = 44890; M = 170; N = 11000; O = 11000; P = 11000; R = "La Ligua";            S = 2200 }
    4 = @{ D = 44890; M = 150; N = 8000;  O = 8000;  P = 8000;  R = "La Ligua";            S = 1600 }
    5 = @{ D = 44890; M = 80;  N = 7000;  O = 7000;  P = 7000;  R = "La Ligua";            S = 1400 }
    6 = @{ D = 44908; M = 110; N = 7000;  O = 7000;  P = 7000;  R = "Provincia de Limarí"; S = 1400 }
    7 = @{ D = 44908; M = 120; N = 6000;  O = 6000;  P = 6000;  R = "Provincia de Limarí"; S = 1200 }
    8 = @{ D = 44908; M = 100; N = 5000;  O = 5000;  P = 5000;  R = "Provincia de Limarí"; S = 1000 }
    9 = @{ D = 44908; M = 120; N = 4000;  O = 4000;  P = 4000;  R = "Provincia de Limarí"; S = 800  }
}

foreach ($row in $newRows.Keys) {
    $vals = $newRows[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
